$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "53.871.31"
$ws.Range("E2").Value = "  +0.58%  "

$ws.Range("D3").Value = "2.246.63"
$ws.Range("E3").Value = "  +2.30%  "

$ws.Range("E4").Value = "  +0.09%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "491.94"
$ws.Range("E5").Value = "  +1.37%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "126.71"
$ws.Range("E6").Value = "  +1.27%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.998"
$ws.Range("E7").Value = "  +0.24%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.524"
$ws.Range("E8").Value = "  +0.68%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.0947"
$ws.Range("E9").Value = "  +3.04%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.151"
$ws.Range("E10").Value = "  +2.43%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.323"
$ws.Range("E11").Value = "  +3.16%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.62"
$ws.Range("E12").Value = "  -0.11%  "

$ws.Range("D13").Value = "2.653.98"
$ws.Range("E13").Value = "  +2.61%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "21.65"
$ws.Range("E14").Value = "  +2.57%  "

$ws.Range("D15").Value = "53.859.08"
$ws.Range("E15").Value = "  +0.69%  "

$ws.Range("E16").Value = "  +0.05%  "

$ws.Range("D17").Value = "2.252.59"
$ws.Range("E17").Value = "  +1.92%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "9.94"
$ws.Range("E18").Value = "  +3.93%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.05"
$ws.Range("E19").Value = "  +2.92%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "298.54"
$ws.Range("E20").Value = "  +1.43%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.39"
$ws.Range("E21").Value = "  +5.24%  "

$ws.Range("E22").Value = "  +0.20%  "

$ws.Range("E23").Value = "  -2.02%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "61.76"
$ws.Range("E24").Value = "  -1.21%  "

$ws.Range("E25").Value = "  +2.21%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.368"
$ws.Range("E26").Value = "  +1.03%  "

$ws.Range("D27").Value = "2.356.60"
$ws.Range("E27").Value = "  +2.47%  "

$ws.Range("E28").Value = "  +1.48%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.00"
$ws.Range("E29").Value = "  -0.06%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "165.38"
$ws.Range("E30").Value = "  +0.19%  "

$ws.Range("E31").Value = "  +0.95%  "

$ws.Range("B32").Value = "USDe"
$ws.Range("C32").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.999"
$ws.Range("E32").Value = "  +0.11%  "

$ws.Range("B33").Value = "PEPE"
$ws.Range("C33").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D33").Value = "0.0₃0672"
$ws.Range("E33").Value = "  +1.27%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.80"
$ws.Range("E34").Value = "  +1.94%  "

$ws.Range("E35").Value = "  +0.19%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.05"
$ws.Range("E36").Value = "  -1.02%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "17.54"
$ws.Range("E37").Value = "  +1.42%  "

$ws.Range("E38").Value = "  +7.14%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.17"
$ws.Range("E39").Value = "  +2.27%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.63"
$ws.Range("E40").Value = "  +2.75%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "35.64"
$ws.Range("E41").Value = "  -0.30%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.38"
$ws.Range("E42").Value = "  +1.87%  "

$ws.Range("E43").Value = "  +0.90%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.33"
$ws.Range("E44").Value = "  +1.96%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.89"
$ws.Range("E45").Value = "  +2.81%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "124.38"
$ws.Range("E46").Value = "  +0.02%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0884"
$ws.Range("E47").Value = "  +1.10%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.537"
$ws.Range("E48").Value = "  +0.78%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "234.17"
$ws.Range("E49").Value = "  +0.83%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0479"
$ws.Range("E50").Value = "  +2.26%  "

$ws.Range("E51").Value = "  +0.04%  "
